$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId=1, the first sheet) ---
$wsExh = $wb.Worksheets.Item("展览")
$wsExh.Range("F2").Value  = 295
$wsExh.Range("F4").Value  = 864
$wsExh.Range("F6").Value  = 318
$wsExh.Range("F7").Value  = 9635
$wsExh.Range("F8").Value  = 81
$wsExh.Range("F11").Value = 122
$wsExh.Range("F17").Value = 278
$wsExh.Range("F18").Value = 787
$wsExh.Range("F20").Value = 95

# --- Sheet "全部类型" (sheetId=4, the fourth sheet) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 295
$wsAll.Range("F4").Value  = 864
$wsAll.Range("F6").Value  = 318
$wsAll.Range("F7").Value  = 9635
$wsAll.Range("F8").Value  = 81
$wsAll.Range("F11").Value = 122
$wsAll.Range("F17").Value = 278
$wsAll.Range("F18").Value = 787
$wsAll.Range("F20").Value = 95
